$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(98, 8).Value = 1605
$ws.Cells.Item(98, 9).Value = 1605
$ws.Cells.Item(98, 10).Value = 0
$ws.Cells.Item(98, 11).Value = 1605
$ws.Cells.Item(98, 12).Value = 0
$ws.Cells.Item(98, 13).Value = -107
$ws.Cells.Item(98, 14).ClearContents()
$ws.Cells.Item(107, 8).Value = 380.5625
$ws.Cells.Item(107, 9).Value = 334.44446
$ws.Cells.Item(107, 10).Value = 439.85715
$ws.Cells.Item(107, 11).Value = 334.44446
$ws.Cells.Item(107, 12).Value = 439.85715
$ws.Cells.Item(107, 13).Value = 1585.55554
$ws.Cells.Item(107, 14).Value = -4279.85715
$ws.Cells.Item(112, 8).Value = 911.5599999999999
$ws.Cells.Item(112, 9).Value = 574.8333
$ws.Cells.Item(112, 11).Value = 1724.4999
$ws.Cells.Item(112, 13).Value = -616.4999
$ws.Cells.Item(122, 8).Value = 1605
$ws.Cells.Item(122, 9).Value = 1605
$ws.Cells.Item(122, 10).Value = 0
$ws.Cells.Item(122, 11).Value = 4815
$ws.Cells.Item(122, 12).Value = 0
$ws.Cells.Item(122, 13).Value = -2365
$ws.Cells.Item(122, 14).ClearContents()
$ws.Cells.Item(127, 8).Value = 795.96295
$ws.Cells.Item(127, 9).Value = 485.33334
$ws.Cells.Item(127, 10).Value = 1417.2222
$ws.Cells.Item(127, 11).Value = 1456.00002
$ws.Cells.Item(127, 12).Value = 4251.6666
$ws.Cells.Item(127, 13).Value = 3503.99998
$ws.Cells.Item(127, 14).Value = -14171.6666
$ws.Cells.Item(129, 8).Value = 6271.3657
$ws.Cells.Item(129, 9).Value = 361.36365
$ws.Cells.Item(129, 10).Value = 8438.366
$ws.Cells.Item(129, 11).Value = 1084.09095
$ws.Cells.Item(129, 12).Value = 25315.098
$ws.Cells.Item(129, 13).Value = 3915.90905
$ws.Cells.Item(129, 14).Value = -35315.098
$ws.Cells.Item(132, 8).Value = 3348.9854
$ws.Cells.Item(132, 9).Value = 3220.2727
$ws.Cells.Item(132, 10).Value = 3893.5386
$ws.Cells.Item(132, 11).Value = 9660.8181
$ws.Cells.Item(132, 12).Value = 11680.6158
$ws.Cells.Item(132, 13).Value = -7130.8181
$ws.Cells.Item(132, 14).Value = -16740.6158
$ws.Cells.Item(138, 8).Value = 3224.169
$ws.Cells.Item(138, 9).Value = 1676.9
$ws.Cells.Item(138, 10).Value = 3830.9412
$ws.Cells.Item(138, 11).Value = 5030.700000000001
$ws.Cells.Item(138, 12).Value = 11492.8236
$ws.Cells.Item(138, 13).Value = 109.2999999999993
$ws.Cells.Item(138, 14).Value = -21772.8236
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 1642.8572
$ws.Cells.Item(2, 10).Value = 1840
$ws.Cells.Item(2, 12).Value = 1840
$ws.Cells.Item(2, 14).Value = -2066
$ws.Cells.Item(32, 8).Value = 10199.75
$ws.Cells.Item(32, 9).Value = 6598.4087
$ws.Cells.Item(32, 10).Value = 19016.828
$ws.Cells.Item(32, 11).Value = 6598.4087
$ws.Cells.Item(32, 12).Value = 19016.828
$ws.Cells.Item(32, 13).Value = -6311.4087
$ws.Cells.Item(32, 14).Value = -19590.828
$ws.Cells.Item(45, 8).Value = 1883.8572
$ws.Cells.Item(45, 9).Value = 2131
$ws.Cells.Item(45, 10).Value = 1266
$ws.Cells.Item(45, 11).Value = 2131
$ws.Cells.Item(45, 12).Value = 1266
$ws.Cells.Item(45, 13).Value = -1754
$ws.Cells.Item(45, 14).Value = -2020
$ws.Cells.Item(61, 8).Value = 1705.5471
$ws.Cells.Item(61, 9).Value = 1478.159
$ws.Cells.Item(61, 10).Value = 2817.2222
$ws.Cells.Item(61, 11).Value = 1478.159
$ws.Cells.Item(61, 12).Value = 2817.2222
$ws.Cells.Item(61, 13).Value = -1266.159
$ws.Cells.Item(61, 14).Value = -3241.2222
$ws.Cells.Item(116, 8).Value = 1642.8572
$ws.Cells.Item(116, 10).Value = 1840
$ws.Cells.Item(116, 12).Value = 1840
$ws.Cells.Item(116, 14).Value = -6428
$ws.Cells.Item(122, 8).Value = 1930.7347
$ws.Cells.Item(122, 9).Value = 1728.125
$ws.Cells.Item(122, 10).Value = 2831.2222
$ws.Cells.Item(122, 11).Value = 5184.375
$ws.Cells.Item(122, 12).Value = 8493.6666
$ws.Cells.Item(122, 13).Value = -2734.375
$ws.Cells.Item(122, 14).Value = -13393.6666
$ws.Cells.Item(132, 8).Value = 19288.824
$ws.Cells.Item(132, 9).Value = 25712.684
$ws.Cells.Item(132, 11).Value = 77138.052
$ws.Cells.Item(132, 13).Value = -74608.052
$ws.Cells.Item(136, 8).Value = 1705.5471
$ws.Cells.Item(136, 9).Value = 1478.159
$ws.Cells.Item(136, 10).Value = 2817.2222
$ws.Cells.Item(136, 11).Value = 4434.477000000001
$ws.Cells.Item(136, 12).Value = 8451.6666
$ws.Cells.Item(136, 13).Value = -1884.477000000001
$ws.Cells.Item(136, 14).Value = -13551.6666
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 1642.8572
$ws.Cells.Item(3, 10).Value = 1840
$ws.Cells.Item(3, 12).Value = 1840
$ws.Cells.Item(3, 14).Value = -2068
$ws.Cells.Item(94, 8).Value = 11857.429
$ws.Cells.Item(94, 9).Value = 7254.75
$ws.Cells.Item(94, 10).Value = 17994.334
$ws.Cells.Item(94, 11).Value = 7254.75
$ws.Cells.Item(94, 12).Value = 17994.334
$ws.Cells.Item(94, 13).Value = -6803.75
$ws.Cells.Item(94, 14).Value = -18896.334
$ws.Cells.Item(105, 8).Value = 2077.389
$ws.Cells.Item(105, 9).Value = 1699.3636
$ws.Cells.Item(105, 11).Value = 1699.3636
$ws.Cells.Item(105, 13).Value = 47.63640000000009
$ws.Cells.Item(107, 8).Value = 14647.5
$ws.Cells.Item(107, 9).Value = 16171.429
$ws.Cells.Item(107, 10).Value = 3980
$ws.Cells.Item(107, 11).Value = 16171.429
$ws.Cells.Item(107, 12).Value = 3980
$ws.Cells.Item(107, 13).Value = -14251.429
$ws.Cells.Item(107, 14).Value = -7820
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(107, 8).Value = 1323.8667
$ws.Cells.Item(107, 9).Value = 1635
$ws.Cells.Item(107, 10).Value = 1051.625
$ws.Cells.Item(107, 11).Value = 1635
$ws.Cells.Item(107, 12).Value = 1051.625
$ws.Cells.Item(107, 13).Value = 285
$ws.Cells.Item(107, 14).Value = -4891.625
$ws.Cells.Item(122, 8).Value = 1142.2916
$ws.Cells.Item(122, 9).Value = 925.1667
$ws.Cells.Item(122, 10).Value = 1359.4166
$ws.Cells.Item(122, 11).Value = 2775.5001
$ws.Cells.Item(122, 12).Value = 4078.2498
$ws.Cells.Item(122, 13).Value = -325.5001000000002
$ws.Cells.Item(122, 14).Value = -8978.2498
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(122, 8).Value = 812.9091
$ws.Cells.Item(122, 10).Value = 1292.7142
$ws.Cells.Item(122, 12).Value = 11634.4278
$ws.Cells.Item(122, 14).Value = -16534.4278
$ws.Cells.Item(131, 8).Value = 1667322
$ws.Cells.Item(131, 9).Value = 3333666.8
$ws.Cells.Item(131, 10).Value = 977.05
$ws.Cells.Item(131, 11).Value = 10001000.4
$ws.Cells.Item(131, 12).Value = 2931.15
$ws.Cells.Item(131, 13).Value = -9995960.399999999
$ws.Cells.Item(131, 14).Value = -13011.15
$ws.Cells.Item(137, 8).Value = 3480.276
$ws.Cells.Item(137, 9).Value = 1133.6923
$ws.Cells.Item(137, 10).Value = 5386.875
$ws.Cells.Item(137, 11).Value = 3401.0769
$ws.Cells.Item(137, 12).Value = 16160.625
$ws.Cells.Item(137, 13).Value = 1698.9231
$ws.Cells.Item(137, 14).Value = -26360.625
$ws.Cells.Item(140, 8).Value = 1240.45
$ws.Cells.Item(140, 9).Value = 831.8125
$ws.Cells.Item(140, 10).Value = 2875
$ws.Cells.Item(140, 11).Value = 2495.4375
$ws.Cells.Item(140, 12).Value = 8625
$ws.Cells.Item(140, 13).Value = 2684.5625
$ws.Cells.Item(140, 14).Value = -18985
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(3, 8).Value = 1638.25
$ws.Cells.Item(3, 9).Value = 3602
$ws.Cells.Item(3, 10).Value = 460
$ws.Cells.Item(3, 11).Value = 3602
$ws.Cells.Item(3, 12).Value = 460
$ws.Cells.Item(3, 13).Value = -3486
$ws.Cells.Item(3, 14).Value = -692
$ws.Cells.Item(122, 8).Value = 2049.4
$ws.Cells.Item(122, 9).Value = 1962.8572
$ws.Cells.Item(122, 10).Value = 2251.3333
$ws.Cells.Item(122, 11).Value = 5888.571599999999
$ws.Cells.Item(122, 12).Value = 6753.999899999999
$ws.Cells.Item(122, 13).Value = -3438.571599999999
$ws.Cells.Item(122, 14).Value = -11653.9999
$ws.Cells.Item(132, 8).Value = 1953.2245
$ws.Cells.Item(132, 9).Value = 1254.0625
$ws.Cells.Item(132, 10).Value = 3269.2942
$ws.Cells.Item(132, 11).Value = 3762.1875
$ws.Cells.Item(132, 12).Value = 9807.882599999999
$ws.Cells.Item(132, 13).Value = -1232.1875
$ws.Cells.Item(132, 14).Value = -14867.8826
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 2159.2
$ws.Cells.Item(7, 9).Value = 1734.8182
$ws.Cells.Item(7, 10).Value = 3326.25
$ws.Cells.Item(7, 11).Value = 1734.8182
$ws.Cells.Item(7, 12).Value = 3326.25
$ws.Cells.Item(7, 13).Value = -1622.8182
$ws.Cells.Item(7, 14).Value = -3550.25
$ws.Cells.Item(61, 8).Value = 1930.5
$ws.Cells.Item(61, 9).Value = 1671.5
$ws.Cells.Item(61, 10).Value = 4002.5
$ws.Cells.Item(61, 11).Value = 1671.5
$ws.Cells.Item(61, 12).Value = 4002.5
$ws.Cells.Item(61, 13).Value = -1469.5
$ws.Cells.Item(61, 14).Value = -4406.5
$ws.Cells.Item(100, 8).Value = 1570.5294
$ws.Cells.Item(100, 9).Value = 1462.5
$ws.Cells.Item(100, 10).Value = 1666.5555
$ws.Cells.Item(100, 11).Value = 1462.5
$ws.Cells.Item(100, 12).Value = 1666.5555
$ws.Cells.Item(100, 13).Value = -921.5
$ws.Cells.Item(100, 14).Value = -2748.5555
$ws.Cells.Item(113, 8).Value = 1930.5
$ws.Cells.Item(113, 9).Value = 1671.5
$ws.Cells.Item(113, 10).Value = 4002.5
$ws.Cells.Item(113, 11).Value = 1671.5
$ws.Cells.Item(113, 12).Value = 4002.5
$ws.Cells.Item(113, 13).Value = 498.5
$ws.Cells.Item(113, 14).Value = -8342.5
$ws.Cells.Item(126, 8).Value = 2159.2
$ws.Cells.Item(126, 9).Value = 1734.8182
$ws.Cells.Item(126, 10).Value = 3326.25
$ws.Cells.Item(126, 11).Value = 5204.4546
$ws.Cells.Item(126, 12).Value = 9978.75
$ws.Cells.Item(126, 13).Value = -2734.4546
$ws.Cells.Item(126, 14).Value = -14918.75
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(25, 8).Value = 0
$ws.Cells.Item(25, 10).Value = 0
$ws.Cells.Item(25, 12).Value = 0
$ws.Cells.Item(25, 14).ClearContents()
$ws.Cells.Item(81, 8).Value = 83336280
$ws.Cells.Item(81, 10).Value = 2611.111
$ws.Cells.Item(81, 12).Value = 5222.222
$ws.Cells.Item(81, 14).Value = -7344.222
$ws.Cells.Item(84, 8).Value = 83336280
$ws.Cells.Item(84, 10).Value = 2611.111
$ws.Cells.Item(84, 12).Value = 26111.11
$ws.Cells.Item(84, 14).Value = -36719.11
$ws.Cells.Item(100, 8).Value = 20599.4
$ws.Cells.Item(100, 9).Value = 50148.75
$ws.Cells.Item(100, 10).Value = 899.8333
$ws.Cells.Item(100, 11).Value = 100297.5
$ws.Cells.Item(100, 12).Value = 1799.6666
$ws.Cells.Item(100, 13).Value = -99756.5
$ws.Cells.Item(100, 14).Value = -2881.6666
$ws.Cells.Item(122, 8).Value = 1064.091
$ws.Cells.Item(122, 9).Value = 883.75
$ws.Cells.Item(122, 11).Value = 2651.25
$ws.Cells.Item(122, 13).Value = -201.25
$ws.Cells.Item(126, 8).Value = 8628
$ws.Cells.Item(126, 9).Value = 8628
$ws.Cells.Item(126, 11).Value = 25884
$ws.Cells.Item(126, 13).Value = -23414
